# Apply the 25-07-2023 data refresh to the SDG national indicators sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the "Last update" timestamp footer (cell A300) ---
$ws.Range("A300").Value = "Last update: 25-07-2023, 10:38"

# --- 2. Fill in the newly published 2022 data point (column S) for the rows that received it ---
#     New values are pasted with column R's number format so the resulting style matches exactly.
$sUpdates = @(
    @{Row=9; Value=412.4}
    @{Row=10; Value=469.6}
    @{Row=11; Value=328.2}
    @{Row=68; Value=36.1}
    @{Row=82; Value=7.6}
    @{Row=83; Value=11.3}
    @{Row=84; Value=9.4}
    @{Row=85; Value=6.4}
    @{Row=86; Value=9}
    @{Row=87; Value=6.6}
    @{Row=88; Value=2.9}
    @{Row=89; Value=3.3}
    @{Row=129; Value=71.900000000000006}
    @{Row=130; Value=77.3}
    @{Row=131; Value=65.099999999999994}
    @{Row=132; Value=859996.9}
    @{Row=146; Value=56.3}
    @{Row=147; Value=64}
    @{Row=148; Value=49.3}
    @{Row=149; Value=27.8}
    @{Row=150; Value=71.3}
    @{Row=151; Value=77.8}
    @{Row=152; Value=51}
    @{Row=153; Value=76.7}
    @{Row=154; Value=35.1}
    @{Row=155; Value=56.4}
    @{Row=156; Value=12.3}
    @{Row=157; Value=56.3}
    @{Row=158; Value=56.3}
    @{Row=159; Value=80.3}
    @{Row=160; Value=59.9}
    @{Row=161; Value=52.5}
    @{Row=162; Value=53.2}
    @{Row=163; Value=16.3}
    @{Row=167; Value=58}
    @{Row=168; Value=50.8}
    @{Row=169; Value=65.8}
    @{Row=170; Value=57.9}
    @{Row=171; Value=58.1}
    @{Row=172; Value=31.1}
    @{Row=173; Value=80.2}
    @{Row=174; Value=35.799999999999997}
    @{Row=175; Value=20.100000000000001}
    @{Row=176; Value=22.1}
    @{Row=177; Value=21.3}
    @{Row=178; Value=22.8}
    @{Row=179; Value=21}
    @{Row=180; Value=23.6}
    @{Row=181; Value=8}
    @{Row=182; Value=8.6}
    @{Row=183; Value=7.5}
    @{Row=184; Value=1939}
    @{Row=185; Value=2234}
    @{Row=186; Value=-295}
    @{Row=207; Value=30.4}
    @{Row=225; Value=6.6}
    @{Row=246; Value=859996.9}
    @{Row=258; Value=0.19800000000000001}
    @{Row=259; Value=31}
    @{Row=289; Value=4}
    @{Row=292; Value=0.53}
    @{Row=294; Value=9.59}
)

foreach ($u in $sUpdates) {
    $r = $u.Row
    $ws.Range("R$r").Copy() | Out-Null
    $ws.Range("S$r").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Range("S$r").Value = $u.Value
}
$excel.CutCopyMode = 0

# --- 3. Row heights were recalculated by Excel on save (autofit reflow) ---
#     Rows 3, 9, 10 and 11 lost their explicit height altogether in the target
#     file (they now inherit the sheet's new default row height of 15), so we
#     line them up with that same effective height here.
$ws.Rows.Item(3).RowHeight = 15
$ws.Rows.Item(9).RowHeight = 15
$ws.Rows.Item(10).RowHeight = 15
$ws.Rows.Item(11).RowHeight = 15

$heightUpdates = @(
    @{Row=4; Height=27}
    @{Row=5; Height=27}
    @{Row=6; Height=27}
    @{Row=7; Height=27}
    @{Row=8; Height=27}
    @{Row=12; Height=36}
    @{Row=13; Height=27}
    @{Row=14; Height=27}
    @{Row=15; Height=27}
    @{Row=16; Height=27}
    @{Row=17; Height=36}
    @{Row=18; Height=18}
    @{Row=19; Height=18}
    @{Row=20; Height=18}
    @{Row=21; Height=18}
    @{Row=22; Height=18}
    @{Row=23; Height=18}
    @{Row=24; Height=18}
    @{Row=25; Height=18}
    @{Row=26; Height=18}
    @{Row=27; Height=18}
    @{Row=28; Height=18}
    @{Row=29; Height=18}
    @{Row=30; Height=18}
    @{Row=31; Height=18}
    @{Row=32; Height=18}
    @{Row=33; Height=18}
    @{Row=34; Height=18}
    @{Row=35; Height=18}
    @{Row=36; Height=18}
    @{Row=37; Height=18}
    @{Row=38; Height=18}
    @{Row=39; Height=18}
    @{Row=40; Height=18}
    @{Row=41; Height=18}
    @{Row=42; Height=18}
    @{Row=43; Height=18}
    @{Row=44; Height=18}
    @{Row=45; Height=18}
    @{Row=46; Height=18}
    @{Row=47; Height=18}
    @{Row=48; Height=18}
    @{Row=49; Height=18}
    @{Row=50; Height=18}
    @{Row=51; Height=18}
    @{Row=52; Height=18}
    @{Row=53; Height=18}
    @{Row=54; Height=18}
    @{Row=55; Height=18}
    @{Row=56; Height=18}
    @{Row=57; Height=18}
    @{Row=58; Height=18}
    @{Row=59; Height=18}
    @{Row=60; Height=18}
    @{Row=61; Height=18}
    @{Row=62; Height=27}
    @{Row=63; Height=27}
    @{Row=64; Height=27}
    @{Row=65; Height=27}
    @{Row=66; Height=18}
    @{Row=67; Height=45}
    @{Row=68; Height=36}
    @{Row=69; Height=36}
    @{Row=70; Height=45}
    @{Row=71; Height=45}
    @{Row=72; Height=45}
    @{Row=73; Height=45}
    @{Row=74; Height=45}
    @{Row=75; Height=45}
    @{Row=76; Height=27}
    @{Row=77; Height=27}
    @{Row=78; Height=27}
    @{Row=79; Height=36}
    @{Row=80; Height=36}
    @{Row=81; Height=36}
    @{Row=82; Height=36}
    @{Row=83; Height=36}
    @{Row=84; Height=36}
    @{Row=85; Height=36}
    @{Row=86; Height=36}
    @{Row=87; Height=36}
    @{Row=88; Height=36}
    @{Row=89; Height=36}
    @{Row=90; Height=27}
    @{Row=91; Height=27}
    @{Row=92; Height=27}
    @{Row=93; Height=27}
    @{Row=94; Height=27}
    @{Row=95; Height=27}
    @{Row=96; Height=27}
    @{Row=97; Height=18}
    @{Row=98; Height=18}
    @{Row=99; Height=18}
    @{Row=100; Height=18}
    @{Row=101; Height=18}
    @{Row=102; Height=18}
    @{Row=103; Height=18}
    @{Row=104; Height=18}
    @{Row=105; Height=18}
    @{Row=106; Height=18}
    @{Row=107; Height=18}
    @{Row=108; Height=18}
    @{Row=109; Height=18}
    @{Row=110; Height=18}
    @{Row=111; Height=18}
    @{Row=112; Height=18}
    @{Row=113; Height=18}
    @{Row=114; Height=18}
    @{Row=115; Height=18}
    @{Row=116; Height=18}
    @{Row=117; Height=18}
    @{Row=118; Height=18}
    @{Row=119; Height=18}
    @{Row=120; Height=36}
    @{Row=121; Height=18}
    @{Row=122; Height=27}
    @{Row=123; Height=27}
    @{Row=124; Height=36}
    @{Row=125; Height=36}
    @{Row=126; Height=18}
    @{Row=127; Height=18}
    @{Row=128; Height=18}
    @{Row=129; Height=27}
    @{Row=130; Height=27}
    @{Row=131; Height=27}
    @{Row=132; Height=18}
    @{Row=133; Height=18}
    @{Row=134; Height=36}
    @{Row=135; Height=36}
    @{Row=136; Height=36}
    @{Row=137; Height=36}
    @{Row=138; Height=27}
    @{Row=139; Height=27}
    @{Row=140; Height=27}
    @{Row=141; Height=18}
    @{Row=142; Height=36}
    @{Row=143; Height=36}
    @{Row=144; Height=36}
    @{Row=145; Height=27}
    @{Row=146; Height=18}
    @{Row=147; Height=18}
    @{Row=148; Height=18}
    @{Row=149; Height=18}
    @{Row=150; Height=18}
    @{Row=151; Height=18}
    @{Row=152; Height=18}
    @{Row=153; Height=18}
    @{Row=154; Height=18}
    @{Row=155; Height=18}
    @{Row=156; Height=18}
    @{Row=157; Height=18}
    @{Row=158; Height=18}
    @{Row=159; Height=18}
    @{Row=160; Height=27}
    @{Row=161; Height=18}
    @{Row=162; Height=18}
    @{Row=163; Height=36}
    @{Row=164; Height=27}
    @{Row=165; Height=27}
    @{Row=166; Height=27}
    @{Row=167; Height=27}
    @{Row=168; Height=27}
    @{Row=169; Height=27}
    @{Row=170; Height=27}
    @{Row=171; Height=27}
    @{Row=172; Height=27}
    @{Row=173; Height=27}
    @{Row=174; Height=27}
    @{Row=175; Height=27}
    @{Row=176; Height=27}
    @{Row=177; Height=27}
    @{Row=178; Height=27}
    @{Row=179; Height=27}
    @{Row=180; Height=27}
    @{Row=181; Height=27}
    @{Row=182; Height=27}
    @{Row=183; Height=27}
    @{Row=184; Height=18}
    @{Row=185; Height=18}
    @{Row=186; Height=18}
    @{Row=187; Height=27}
    @{Row=188; Height=27}
    @{Row=189; Height=27}
    @{Row=190; Height=27}
    @{Row=191; Height=27}
    @{Row=192; Height=18}
    @{Row=193; Height=18}
    @{Row=194; Height=18}
    @{Row=195; Height=18}
    @{Row=196; Height=18}
    @{Row=197; Height=18}
    @{Row=198; Height=18}
    @{Row=199; Height=18}
    @{Row=200; Height=18}
    @{Row=201; Height=18}
    @{Row=202; Height=18}
    @{Row=203; Height=18}
    @{Row=204; Height=27}
    @{Row=205; Height=27}
    @{Row=206; Height=54}
    @{Row=207; Height=54}
    @{Row=208; Height=54}
    @{Row=209; Height=54}
    @{Row=210; Height=54}
    @{Row=211; Height=54}
    @{Row=212; Height=54}
    @{Row=213; Height=54}
    @{Row=214; Height=54}
    @{Row=215; Height=36}
    @{Row=216; Height=36}
    @{Row=217; Height=36}
    @{Row=218; Height=36}
    @{Row=219; Height=45}
    @{Row=220; Height=45}
    @{Row=221; Height=45}
    @{Row=222; Height=27}
    @{Row=223; Height=27}
    @{Row=224; Height=27}
    @{Row=225; Height=27}
    @{Row=226; Height=27}
    @{Row=227; Height=27}
    @{Row=228; Height=27}
    @{Row=229; Height=27}
    @{Row=230; Height=27}
    @{Row=231; Height=36}
    @{Row=232; Height=36}
    @{Row=233; Height=36}
    @{Row=234; Height=36}
    @{Row=235; Height=27}
    @{Row=236; Height=36}
    @{Row=237; Height=36}
    @{Row=238; Height=27}
    @{Row=239; Height=27}
    @{Row=240; Height=27}
    @{Row=241; Height=27}
    @{Row=242; Height=27}
    @{Row=243; Height=27}
    @{Row=244; Height=27}
    @{Row=245; Height=18}
    @{Row=246; Height=18}
    @{Row=247; Height=18}
    @{Row=248; Height=18}
    @{Row=249; Height=135}
    @{Row=250; Height=36}
    @{Row=251; Height=72}
    @{Row=252; Height=72}
    @{Row=253; Height=72}
    @{Row=254; Height=72}
    @{Row=255; Height=18}
    @{Row=256; Height=27}
    @{Row=257; Height=36}
    @{Row=258; Height=36}
    @{Row=259; Height=36}
    @{Row=260; Height=36}
    @{Row=261; Height=36}
    @{Row=262; Height=36}
    @{Row=263; Height=36}
    @{Row=264; Height=36}
    @{Row=265; Height=27}
    @{Row=266; Height=27}
    @{Row=267; Height=27}
    @{Row=268; Height=27}
    @{Row=269; Height=27}
    @{Row=270; Height=27}
    @{Row=271; Height=27}
    @{Row=272; Height=27}
    @{Row=273; Height=27}
    @{Row=274; Height=27}
    @{Row=275; Height=27}
    @{Row=276; Height=27}
    @{Row=277; Height=27}
    @{Row=278; Height=27}
    @{Row=279; Height=27}
    @{Row=280; Height=27}
    @{Row=281; Height=27}
    @{Row=282; Height=27}
    @{Row=283; Height=27}
    @{Row=284; Height=27}
    @{Row=285; Height=27}
    @{Row=286; Height=27}
    @{Row=287; Height=27}
    @{Row=288; Height=27}
    @{Row=289; Height=27}
    @{Row=290; Height=36}
    @{Row=291; Height=36}
    @{Row=292; Height=72}
    @{Row=293; Height=27}
    @{Row=294; Height=27}
    @{Row=295; Height=27}
    @{Row=296; Height=27}
    @{Row=297; Height=27}
)

foreach ($u in $heightUpdates) {
    $ws.Rows.Item($u.Row).RowHeight = $u.Height
}

# --- 4. Selection / active cell moved from A1:B1 to A4 ---
$ws.Activate()
$ws.Range("A4").Select() | Out-Null

# --- 5. Print scale bumped from 63% to 65% ---
$ws.PageSetup.Zoom = 65

Write-Output "edit complete"
